$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1423809062965233
$ws.Range("D2").Value = 0.2851467616826824
$ws.Range("E2").Value = 0.1372607365515108
$ws.Range("F2").Value = 6.825672395130994
$ws.Range("G2").Value = 0.002695213221253804
$ws.Range("J2").Value = 0.1367055720989043
$ws.Range("M2").Value = 6.549288000167593
$ws.Range("B3").Value = 0.1329604150238453
$ws.Range("D3").Value = 0.2560228835505711
$ws.Range("E3").Value = 0.1193733527267895
$ws.Range("F3").Value = 6.847732267853075
$ws.Range("G3").Value = 0.00271376776157433
$ws.Range("J3").Value = 0.1264069728587884
$ws.Range("M3").Value = 5.964698587850791
$ws.Range("B4").Value = 0.1272503147543489
$ws.Range("D4").Value = 0.2385168469926953
$ws.Range("E4").Value = 0.1084035697961028
$ws.Range("F4").Value = 6.872095838575831
$ws.Range("G4").Value = 0.002725658888421509
$ws.Range("J4").Value = 0.1201531893065777
$ws.Range("M4").Value = 5.607343347165113
$ws.Range("B5").Value = 0.1249420585963747
$ws.Range("D5").Value = 0.2314715214307626
$ws.Range("E5").Value = 0.1039352463805727
$ws.Range("F5").Value = 6.884682676537011
$ws.Range("G5").Value = 0.002730631112166133
$ws.Range("J5").Value = 0.1176215346582268
$ws.Range("M5").Value = 5.462076064808059
$ws.Range("B6").Value = 0.1245599018989338
$ws.Range("D6").Value = 0.2303068244362976
$ws.Range("E6").Value = 0.1031933644713874
$ws.Range("F6").Value = 6.886931505124608
$ws.Range("G6").Value = 0.002731464418329869
$ws.Range("J6").Value = 0.117202151867545
$ws.Range("M6").Value = 5.437975018746727
$ws.Range("B7").Value = 0.1272191092360941
$ws.Range("D7").Value = 0.2384214806704961
$ws.Range("E7").Value = 0.1083433022539495
$ws.Range("F7").Value = 6.872254905279163
$ws.Range("G7").Value = 0.002725725431888926
$ws.Range("J7").Value = 0.1201189793161603
$ws.Range("M7").Value = 5.605382824650832
$ws.Range("B8").Value = 0.139117342583873
$ws.Range("D8").Value = 0.2750230562874947
$ws.Range("E8").Value = 0.1310895965230898
$ws.Range("F8").Value = 6.830995242097657
$ws.Range("G8").Value = 0.002701508037627576
$ws.Range("J8").Value = 0.1331397333637625
$ws.Range("M8").Value = 6.347367134067355
$ws.Range("B9").Value = 0.1630388474108599
$ws.Range("D9").Value = 0.3500705838918918
$ws.Range("E9").Value = 0.1758660676599249
$ws.Range("F9").Value = 6.838833463057796
$ws.Range("G9").Value = 0.00265792094052208
$ws.Range("J9").Value = 0.1592612256124397
$ws.Range("M9").Value = 7.81700964078118
$ws.Range("B10").Value = 0.1809771585331532
$ws.Range("D10").Value = 0.4076235896185665
$ws.Range("E10").Value = 0.2089689795170386
$ws.Range("F10").Value = 6.902906266646369
$ws.Range("G10").Value = 0.002628200780592222
$ws.Range("J10").Value = 0.1788673011249813
$ws.Range("M10").Value = 8.908795825016455
$ws.Range("B11").Value = 0.1892178160615998
$ws.Range("D11").Value = 0.4344284048614782
$ws.Range("E11").Value = 0.22409714357984
$ws.Range("F11").Value = 6.945692305268949
$ws.Range("G11").Value = 0.002615163153775613
$ws.Range("J11").Value = 0.1878905501042851
$ws.Range("M11").Value = 9.408845529922132
$ws.Range("B12").Value = 0.1923499719273991
$ws.Range("D12").Value = 0.4446768717425016
$ws.Range("E12").Value = 0.2298378581183442
$ws.Range("F12").Value = 6.963940893049084
$ws.Range("G12").Value = 0.00261029403629461
$ws.Range("J12").Value = 0.1913236488665149
$ws.Range("M12").Value = 9.59875412065287
$ws.Range("B13").Value = 0.1916748896785947
$ws.Range("D13").Value = 0.4424651810825821
$ws.Range("E13").Value = 0.2286009249754244
$ws.Range("F13").Value = 6.959918312059301
$ws.Range("G13").Value = 0.00261133968852439
$ws.Range("J13").Value = 0.1905835296212928
$ws.Range("M13").Value = 9.557828366510989
$ws.Range("B14").Value = 0.189475267750538
$ws.Range("D14").Value = 0.435269537351985
$ws.Range("E14").Value = 0.2245691829558325
$ws.Range("F14").Value = 6.94715210970719
$ws.Range("G14").Value = 0.002614761214889408
$ws.Range("J14").Value = 0.1881726615865489
$ws.Range("M14").Value = 9.424458011521608
$ws.Range("B15").Value = 0.1881294473117094
$ws.Range("D15").Value = 0.4308750254926395
$ws.Range("E15").Value = 0.2221012503414528
$ws.Range("F15").Value = 6.939601554850753
$ws.Range("G15").Value = 0.002616865805180688
$ws.Range("J15").Value = 0.1866980806164804
$ws.Range("M15").Value = 9.342838473127642
$ws.Range("B16").Value = 0.1804402331171389
$ws.Range("D16").Value = 0.4058850877510167
$ws.Range("E16").Value = 0.2079818858674685
$ws.Range("F16").Value = 6.900391728517491
$ws.Range("G16").Value = 0.002629062405810159
$ws.Range("J16").Value = 0.1782798093836675
$ws.Range("M16").Value = 8.876189917361899
$ws.Range("B17").Value = 0.1757437685213432
$ws.Range("D17").Value = 0.3907201265220408
$ws.Range("E17").Value = 0.1993393018991299
$ws.Range("F17").Value = 6.8798902433102
$ws.Range("G17").Value = 0.00263666715974589
$ws.Range("J17").Value = 0.173143010780592
$ws.Range("M17").Value = 8.59082803399707
$ws.Range("B18").Value = 0.1730500541436015
$ws.Range("D18").Value = 0.3820558156992035
$ws.Range("E18").Value = 0.1943747278377259
$ws.Range("F18").Value = 6.869376741131504
$ws.Range("G18").Value = 0.002641086689023415
$ws.Range("J18").Value = 0.1701982048902693
$ws.Range("M18").Value = 8.427010536455384
$ws.Range("B19").Value = 0.1721393077329054
$ws.Range("D19").Value = 0.3791319723541164
$ws.Range("E19").Value = 0.1926948503139698
$ws.Range("F19").Value = 6.866034095536662
$ws.Range("G19").Value = 0.002642590917923293
$ws.Range("J19").Value = 0.1692027833804843
$ws.Range("M19").Value = 8.371597018793295
$ws.Range("B20").Value = 0.1762429317549419
$ws.Range("D20").Value = 0.3923283825089356
$ws.Range("E20").Value = 0.2002586419577881
$ws.Range("F20").Value = 6.881939748187847
$ws.Range("G20").Value = 0.002635852924138223
$ws.Range("J20").Value = 0.173688814938771
$ws.Range("M20").Value = 8.621172238483211
$ws.Range("B21").Value = 0.1901210349619333
$ws.Range("D21").Value = 0.4373803417584838
$ws.Range("E21").Value = 0.2257530607176506
$ws.Range("F21").Value = 6.950845611505315
$ws.Range("G21").Value = 0.002613754396980925
$ws.Range("J21").Value = 0.1888803426313217
$ws.Range("M21").Value = 9.46361664917157
$ws.Range("B22").Value = 0.199258800955846
$ws.Range("D22").Value = 0.4674002305033014
$ws.Range("E22").Value = 0.2424861569064234
$ws.Range("F22").Value = 7.00785067026311
$ws.Range("G22").Value = 0.00259970710014204
$ws.Range("J22").Value = 0.1989038751040795
$ws.Range("M22").Value = 10.01744474609592
$ws.Range("B23").Value = 0.1943756022120482
$ws.Range("D23").Value = 0.4513224696289058
$ws.Range("E23").Value = 0.2335481889279549
$ws.Range("F23").Value = 6.976301296694373
$ws.Range("G23").Value = 0.002607168706565174
$ws.Range("J23").Value = 0.1935450177664393
$ws.Range("M23").Value = 9.721537607064192
$ws.Range("B24").Value = 0.1760172402937314
$ws.Range("D24").Value = 0.3916011216998356
$ws.Range("E24").Value = 0.1998429955024648
$ws.Range("F24").Value = 6.88100920807824
$ws.Range("G24").Value = 0.002636220892217716
$ws.Range("J24").Value = 0.1734420307779487
$ws.Range("M24").Value = 8.607452878152287
$ws.Range("B25").Value = 0.1565039960823782
$ws.Range("D25").Value = 0.329375736453585
$ws.Range("E25").Value = 0.1637244480584812
$ws.Range("F25").Value = 6.826789048409353
$ws.Range("G25").Value = 0.002669302260871575
$ws.Range("J25").Value = 0.1521262273806627
$ws.Range("M25").Value = 5.462076064808059
